# add reservation status in create/update booking
#
# 1) Shrinks the sheet's small Arial header/data font from 6pt to 5pt
#    (it now has one more column to share the row with).
# 2) Adds a new "Resv Status" column (Q) with value "RESERVED" on the
#    single booking row.
# 3) Re-tightens all the column widths to fit the wider, 17-column layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Font size 6 -> 5 on every populated cell that uses the small Arial font.
#     (Done in same-style contiguous blocks so the untouched fill-only cells
#     B2/J2/L2/O2 are left exactly as they were.)
$ws.Range("A1:P1").Font.Size = 5
$ws.Range("A2").Font.Size = 5
$ws.Range("C2:I2").Font.Size = 5
$ws.Range("K2").Font.Size = 5
$ws.Range("M2:N2").Font.Size = 5
$ws.Range("P2").Font.Size = 5

# --- 2) New column Q: "Resv Status" header + "RESERVED" value.
#     Copy the "Room Tax" cell's look first so Q1/Q2 pick up the same
#     header/data cell style, then overwrite the copied text.
$ws.Range("P1").Copy($ws.Range("Q1"))
$ws.Range("Q1").Value = "Resv Status"

$ws.Range("P2").Copy($ws.Range("Q2"))
$ws.Range("Q2").Value = "RESERVED"

# --- 3) Re-fit column widths for the new 17-column layout.
$ws.Columns.Item(1).ColumnWidth = 8
$ws.Columns.Item(2).ColumnWidth = 12.67
$ws.Columns.Item(3).ColumnWidth = 8.67
$ws.Columns.Item(4).ColumnWidth = 11
$ws.Columns.Item(5).ColumnWidth = 12
$ws.Columns.Item(6).ColumnWidth = 6.67
$ws.Columns.Item(7).ColumnWidth = 7.67
$ws.Columns.Item(8).ColumnWidth = 12
$ws.Columns.Item(9).ColumnWidth = 8
$ws.Columns.Item(10).ColumnWidth = 5.33
$ws.Columns.Item(11).ColumnWidth = 4.67
$ws.Columns.Item(12).ColumnWidth = 9
$ws.Columns.Item(13).ColumnWidth = 12
$ws.Columns.Item(14).ColumnWidth = 9
$ws.Columns.Item(15).ColumnWidth = 11.17
$ws.Columns.Item(16).ColumnWidth = 6.83
$ws.Columns.Item(17).ColumnWidth = 9.67
